$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 161.279784
$ws.Range("H2").Value = 483.839352
$ws.Range("I2").Value = 0.3023989599621841
$ws.Range("J2").Value = 0.3023989599621841
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 57.77686466666668
$ws.Range("N2").Value = 173.330594
$ws.Range("O2").Value = 0.5795494107546101
$ws.Range("P2").Value = 0.57954941075461
$ws.Range("Q2").Value = 9318.240253637232
$ws.Range("R2").Value = 83864.1622827351
$ws.Range("S2").Value = 0.1752551390588907
$ws.Range("T2").Value = 0.1752551390588907

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 161.279784
$ws.Range("H3").Value = 483.839352
$ws.Range("I3").Value = 0.3023989599621841
$ws.Range("J3").Value = 0.3023989599621841
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.278736333333333
$ws.Range("N3").Value = 27.836209
$ws.Range("O3").Value = 0.09307334701450438
$ws.Range("P3").Value = 0.09307334701450438
$ws.Range("Q3").Value = 1496.472591632952
$ws.Range("R3").Value = 13468.25332469657
$ws.Range("S3").Value = 0.02814528333738558
$ws.Range("T3").Value = 0.02814528333738558

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 161.279784
$ws.Range("H4").Value = 483.839352
$ws.Range("I4").Value = 0.3023989599621841
$ws.Range("J4").Value = 0.3023989599621841
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 32.637132
$ws.Range("N4").Value = 97.91139600000001
$ws.Range("O4").Value = 0.3273772422308855
$ws.Range("P4").Value = 0.3273772422308855
$ws.Range("Q4").Value = 5263.709599339488
$ws.Range("R4").Value = 47373.38639405539
$ws.Range("S4").Value = 0.09899853756590779
$ws.Range("T4").Value = 0.09899853756590779

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 288.7700093333333
$ws.Range("H5").Value = 866.3100279999999
$ws.Range("I5").Value = 0.541442630470476
$ws.Range("J5").Value = 0.5414426304704759
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 57.77686466666668
$ws.Range("N5").Value = 173.330594
$ws.Range("O5").Value = 0.5795494107546101
$ws.Range("P5").Value = 0.57954941075461
$ws.Range("Q5").Value = 16684.22574904407
$ws.Range("R5").Value = 150158.0317413966
$ws.Range("S5").Value = 0.3137927574465905
$ws.Range("T5").Value = 0.3137927574465904

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 288.7700093333333
$ws.Range("H6").Value = 866.3100279999999
$ws.Range("I6").Value = 0.541442630470476
$ws.Range("J6").Value = 0.5414426304704759
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 9.278736333333333
$ws.Range("N6").Value = 27.836209
$ws.Range("O6").Value = 0.09307334701450438
$ws.Range("P6").Value = 0.09307334701450438
$ws.Range("Q6").Value = 2679.420777578205
$ws.Range("R6").Value = 24114.78699820385
$ws.Range("S6").Value = 0.05039387783422468
$ws.Range("T6").Value = 0.05039387783422467

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 288.7700093333333
$ws.Range("H7").Value = 866.3100279999999
$ws.Range("I7").Value = 0.541442630470476
$ws.Range("J7").Value = 0.5414426304704759
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 32.637132
$ws.Range("N7").Value = 97.91139600000001
$ws.Range("O7").Value = 0.3273772422308855
$ws.Range("P7").Value = 0.3273772422308855
$ws.Range("Q7").Value = 9424.624912253232
$ws.Range("R7").Value = 84821.62421027909
$ws.Range("S7").Value = 0.1772559951896609
$ws.Range("T7").Value = 0.1772559951896608

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 83.28466000000002
$ws.Range("H8").Value = 249.85398
$ws.Range("I8").Value = 0.15615840956734
$ws.Range("J8").Value = 0.15615840956734
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 57.77686466666668
$ws.Range("N8").Value = 173.330594
$ws.Range("O8").Value = 0.5795494107546101
$ws.Range("P8").Value = 0.57954941075461
$ws.Range("Q8").Value = 4811.926529629349
$ws.Range("R8").Value = 43307.33876666413
$ws.Range("S8").Value = 0.09050151424912893
$ws.Range("T8").Value = 0.09050151424912892

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 83.28466000000002
$ws.Range("H9").Value = 249.85398
$ws.Range("I9").Value = 0.15615840956734
$ws.Range("J9").Value = 0.15615840956734
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 9.278736333333333
$ws.Range("N9").Value = 27.836209
$ws.Range("O9").Value = 0.09307334701450438
$ws.Range("P9").Value = 0.09307334701450438
$ws.Range("Q9").Value = 772.7764007513134
$ws.Range("R9").Value = 6954.987606761821
$ws.Range("S9").Value = 0.01453418584289413
$ws.Range("T9").Value = 0.01453418584289413

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 83.28466000000002
$ws.Range("H10").Value = 249.85398
$ws.Range("I10").Value = 0.15615840956734
$ws.Range("J10").Value = 0.15615840956734
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 32.637132
$ws.Range("N10").Value = 97.91139600000001
$ws.Range("O10").Value = 0.3273772422308855
$ws.Range("P10").Value = 0.3273772422308855
$ws.Range("Q10").Value = 2718.172441995121
$ws.Range("R10").Value = 24463.55197795609
$ws.Range("S10").Value = 0.05112270947531689
$ws.Range("T10").Value = 0.05112270947531689
